# "Lesson 17 page object" - add Sex / Region / City columns (E:G) to the
# Registration sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1, E1:G1): bold + centered, like the existing headers ---
$headerRange = $ws.Range("E1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4108     # xlCenter

$ws.Range("E1").Value = "Sex"
$ws.Range("F1").Value = "Region"
$ws.Range("G1").Value = "City"

# --- Data rows 2-4, columns E-G: left aligned, like the existing data cells ---
# (E2 is set up further below with a left-over date number format instead of
#  the plain left-aligned style)
$ws.Range("E3:E4").HorizontalAlignment = -4131   # xlLeft
$ws.Range("F2:G4").HorizontalAlignment = -4131   # xlLeft

$ws.Range("E2").Value = "Мужской"
$ws.Range("F2").Value = "Калининградская область"
$ws.Range("G2").Value = "Гусев"

$ws.Range("E3").Value = "Женский"
$ws.Range("F3").Value = "Московская область"
$ws.Range("G3").Value = "Домодедово"

$ws.Range("E4").Value = "Мужской"
$ws.Range("F4").Value = "Санкт-Петербург"
$ws.Range("G4").Value = "Санкт-Петербург"

# E2 keeps a leftover short-date number format (numFmtId 14)
$ws.Range("E2").NumberFormat = "mm-dd-yy"

# --- Column widths for the new columns ---
$ws.Range("E1:E4").ColumnWidth = 17.5
$ws.Range("F1:F4").ColumnWidth = 22.83
$ws.Range("G1:G4").ColumnWidth = 14.5

# --- Final selection ---
$ws.Range("H8").Select() | Out-Null
